# "Revert to 2.1.1 files" — re-creates the Texas-specific edits to the
# Discount Rate workbook:
#   * inserts a new "Texas Notes" sheet between "About" and "DR"
#   * updates the DR discount-rate value from 3% to 5.87% (VCE WISdom number)
#   * nudges the saved selections on "About" / "DR" to match the authored file

$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")

# Insert the new "Texas Notes" sheet right after "About" (i.e. before "DR"),
# so the sheet order becomes About, Texas Notes, DR.
$notesSheet = $wb.Worksheets.Add($null, $aboutSheet)
$notesSheet.Name = "Texas Notes"

$notesSheet.Range("A1").Value = "updated to the VCE WISdom number"
$notesSheet.Range("A2").Value = 0.0587
$notesSheet.Range("A4").Value = "their feedback was the 3% was a bit low"
$notesSheet.Range("A5").Select() | Out-Null

# Re-fetch "DR" by name now that the sheet collection has shifted.
$drSheet = $wb.Worksheets.Item("DR")

# Update the discount rate itself.
$drSheet.Range("B2").Value = 0.0587
$drSheet.Range("B2").Select() | Out-Null

$aboutSheet.Range("C23").Select() | Out-Null

Write-Output "Texas Notes sheet inserted; DR discount rate updated to 5.87%"
